# Agregar solucionario del trabajo 1
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Tabla 1")
$ws2 = $wb.Worksheets.Item("Tabla 2")
$ws3 = $wb.Worksheets.Item("Tabla 3")

# Update H column values on "Tabla 1" (truth table solution values)
$ws1.Range("H4").Value = "v"
$ws1.Range("H5").Value = "v"
$ws1.Range("H6").Value = "f"
$ws1.Range("H7").Value = "f"

# Update selections per sheet
[void]$ws1.Range("H8").Select()
[void]$ws2.Range("L13").Select()
[void]$ws3.Range("G30").Select()

# Activate "Tabla 2" as the active sheet (tabSelected) - do this last so its
# selection (L13) is the one left "current" / visible when the workbook is saved
$ws2.Activate()
[void]$ws2.Range("L13").Select()
